$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 51 (shifts existing rows 51.. down by one)
$ws.Rows.Item(51).Insert()

# Copy formatting from the row above (row 50) into the new row 51 so styles
# (e.g. the date format on column D) carry over. Restrict to the used
# columns (A:R) so we don't stamp style info across the entire 16384-column
# row (which would blow up the sheet dimension).
$ws.Range("A50:R50").Copy()
$ws.Range("A51:R51").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's values
$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = 44946
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 100112030
$ws.Range("G51").Value = "Poroto granado"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 50
$ws.Range("K51").Value = 45000
$ws.Range("L51").Value = 45000
$ws.Range("M51").Value = 45000
$ws.Range("N51").Value = "$/saco 25 kilos"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 1800
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
